$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 with same style as the other header cells (A1:C1)
$ws.Range("D1").Value = "Tipo"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update existing MSE/R2 values and add "multiple" type column for each row
$ws.Range("B2").Value = 0.5256126682387602
$ws.Range("C2").Value = 0.9895339355731295
$ws.Range("D2").Value = "multiple"

$ws.Range("B3").Value = 0.235593565623943
$ws.Range("C3").Value = 0.9953953178351463
$ws.Range("D3").Value = "multiple"

$ws.Range("B4").Value = 0.2782160958434385
$ws.Range("C4").Value = 0.9946477409688294
$ws.Range("D4").Value = "multiple"

$ws.Range("B5").Value = 0.41252899801167
$ws.Range("C5").Value = 0.9918654929617656
$ws.Range("D5").Value = "multiple"
